# Update statistics in the "展览" (exhibition) sheet and the combined
# "全部类型" (all types) sheet. Both sheets list the same four exhibition
# events; column F holds a numeric stat (e.g. registrations/views) that
# needs to be refreshed.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" : rows 2-5, column F ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 430
$wsExpo.Range("F3").Value = 2938
$wsExpo.Range("F4").Value = 135
$wsExpo.Range("F5").Value = 67

# --- Sheet "全部类型" : rows 2, 7, 8, 10, column F (the same four events) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 430
$wsAll.Range("F7").Value = 2938
$wsAll.Range("F8").Value = 135
$wsAll.Range("F10").Value = 67
